$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.692.08"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.749.65"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.40%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "612.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.31"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.59%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.743.75"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.57%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.166"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.33"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.73%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.491"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.69%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "40.73"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.18%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000252"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.47%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.374.74"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.743.03"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.60%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.713.70"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.25%  "
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.57"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.78%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "512.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.67"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.51"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.68%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.724"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.45%  "
$ws.Range("B24").Value = "Fetch.AI"
$ws.Range("C24").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.51"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.27%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "87.78"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.82%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.26"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.39%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.21%  "
$ws.Range("E28").Value = "  +17.10%  "
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.48"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.89%  "
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.83"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.61%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.82"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.94%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.36"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.89%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.115"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.41%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.19"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.88%  "
$ws.Range("E37").Value = "  +1.32%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.337"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.88%  "
$ws.Range("E39").Value = "  +3.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.132"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.99%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "51.18"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "44.80"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.48%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.77"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.12%  "
$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "420.29"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.81%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.057.43"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.69%  "
$ws.Range("E46").Value = "  -3.79%  "
$ws.Range("E47").Value = "  -0.44%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "27.75"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.15%  "
$ws.Range("B49").Value = "ThetaToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.50"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.09%  "
$ws.Range("B50").Value = "USDe"
$ws.Range("C50").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "134.91"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.14%  "
